# "fix and up load carousel"
#
# The "carousel" task row (row 4) is updated:
#   - Build start/finish dates become plain text dates (matching the
#     "dd/m/yyyy" style notes used elsewhere in the sheet) instead of
#     real date serials.
#   - The assigned builder/tester names are corrected to
#     "Nguyễn Trí Hậu" / "Nguyễn Văn Lanh".
# The active selection is also left on F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "18/8/2022"
$ws.Range("C4").Value = "19/8/2022"
$ws.Range("E4").Value = "Nguyễn Trí Hậu"
$ws.Range("I4").Value = "Nguyễn Văn Lanh"

[void]$ws.Range("F6").Select()
